# "all done with time measurment on tables wo keys and index"
#
# Fill in row 3 (previously-blank summary/header values) on the
# "Table content drop wo I&Keys" sheet with the measured totals, and
# make that sheet the active tab/selection (it was previously on
# "Table content drop w I&Keys").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table content drop wo I&Keys")

# Row 3 values (A..N) — totals that were missing before.
$ws.Cells.Item(3, 1).Value  = 17449200   # A3
$ws.Cells.Item(3, 2).Value  = 403186     # B3
$ws.Cells.Item(3, 3).Value  = 69409      # C3
$ws.Cells.Item(3, 4).Value  = 83584      # D3
$ws.Cells.Item(3, 5).Value  = 4557053    # E3
$ws.Cells.Item(3, 6).Value  = 104503     # F3
$ws.Cells.Item(3, 7).Value  = 17449200   # G3
$ws.Cells.Item(3, 8).Value  = 1427       # H3
$ws.Cells.Item(3, 9).Value  = 17449200   # I3
$ws.Cells.Item(3, 10).Value = 3413       # J3
$ws.Cells.Item(3, 11).Value = 69409      # K3 (was 79738)
$ws.Cells.Item(3, 12).Value = 18111      # L3 (was 23580)
$ws.Cells.Item(3, 13).Value = 69409      # M3 (was 79738)
$ws.Cells.Item(3, 14).Value = 15218      # N3 (was 20500)

# Make this sheet the active one and move the selection to F3
# (previously "Table content drop w I&Keys" / B12 was active).
[void]$ws.Activate()
[void]$ws.Range("F3").Select()
